# Update Sage scrape results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "The" column (old column G). Deleting the entire column
# shifts H:Z left into G:Y, which is exactly the column 1 header
# realignment and the Y/Z numeric-column collapse shown in the diff.
$ws.Range("G:G").Delete()

# Year values are stored as text in this sheet, so a leading apostrophe
# keeps numeric-looking entries ("2017", "2014", ...) as text instead of
# letting Excel auto-convert them to numbers.

# Row 2
$ws.Range("B2").Value = "Digital Assays Part II: Digital Protein and Cell Assays"
$ws.Range("C2").Value = "Amar S. Basu"
$ws.Range("D2").Value = "'2017"
$ws.Range("E2").Value = "10.1177/2472630317705681"

# Row 3
$ws.Range("B3").Value = "Human Factors in Cyber Warfare II: Emerging Perspectives"
$ws.Range("C3").Value = "Panel Chair: Dr.Vincent F. Mancuso, Panelists: Dr.James C. Christensen, Dr.Jennifer Cowley, Dr.Victor Finomore, Prof.Cleotide Gonzalez, Dr.Benjamin Knott"
$ws.Range("D3").Value = "'2014"
$ws.Range("E3").Value = "10.1177/1541931214581085"

# Row 4
$ws.Range("B4").Value = "The dynamics of cyber conflict between rival antagonists, 2001–11"
$ws.Range("C4").Value = "Brandon Valeriano, Ryan C Maness"
$ws.Range("D4").Value = "'2014"
$ws.Range("E4").Value = "10.1177/0022343313518940"

# Row 5
$ws.Range("B5").Value = "Towards a Chronology of Robotic Art"
$ws.Range("C5").Value = "Eduardo Kac"
$ws.Range("D5").Value = "'2001"
$ws.Range("E5").Value = "10.1177/135485650100700109"
$ws.Range("F5").Value = "Restricted"

# Row 6
$ws.Range("B6").Value = "Warring from the virtual to the real: Assessing the public’s threshold for war over cyber security"
$ws.Range("C6").Value = "Sarah Kreps, Debak Das"
$ws.Range("D6").Value = "'2017"
$ws.Range("E6").Value = "10.1177/2053168017715930"
$ws.Range("F6").Value = "Open Access"

# Row 7
$ws.Range("B7").Value = "An Adversarial Model for Expressing Attacks on Control Protocols"
$ws.Range("C7").Value = "Jonathan Butts, Mason Rice, Sujeet Shenoi"
$ws.Range("D7").Value = "'2012"
$ws.Range("E7").Value = "10.1177/1548512911449409"
$ws.Range("F7").Value = "Restricted"

# Row 8
$ws.Range("B8").Value = "Simplification and Linearization of Manipulator Dynamics by the Design of Inertia Distribution"
$ws.Range("C8").Value = "D.C.H. Yang, S.W. Tzeng"
$ws.Range("D8").Value = "'1986"
$ws.Range("E8").Value = "10.1177/027836498600500307"

# Row 9
$ws.Range("B9").Value = "A novel ensemble learning approach for fault detection of sensor data in cyber-physical system"
$ws.Range("C9").Value = "Ramesh Sneka Nandhini, Ramanathan Lakshmanan"
$ws.Range("E9").Value = "10.3233/JIFS-235809"

# Row 10
$ws.Range("B10").Value = "From the Ontology of Video Games to the Epistemology of Digital Movements. Towards a Semiotics of Virtual Practices"
$ws.Range("C10").Value = "Enzo D’Armenio"
$ws.Range("D10").Value = "'2024"
$ws.Range("E10").Value = "10.1177/15554120241263630"

# Row 11
$ws.Range("B11").Value = "Global versus Local Optimization in Redundancy Resolution of Robotic Manipulators"
$ws.Range("C11").Value = "Kazem Kazerounian, Zhaoyu Wang"
$ws.Range("D11").Value = "'1988"
$ws.Range("E11").Value = "10.1177/027836498800700501"
